$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Goofellow", "Ian", 5, 2),
    @("Kurakin", "Alexei", 5, 2),
    @("Bengio", "Samy", 5, 2),
    @("Madry", "Aleksander", 5, 3),
    @("Papernot", "Nicolas", 5, 3),
    @("Carlini", "Nicholas", 5, 3),
    @("Wagner", "David", 5, 3)
)

$row = 28
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

$ws.Range("B35").Select()
